$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    231 = "At this level the team is not cross-functional in nature. Some of the skills required to complete end-to-end customer features are missing and there is dependence of teams/people outside the team."
    232 = "Team members are not 100% allocated to the teams and have other priorities as well that they have to allocate their time to."
    233 = "Collaboration between team members is low at this mastery level and mainly ad hoc."
    234 = "Software release requires multiple handoffs and dependency on external teams and stakeholders. The team does not know the different external teams and stake holders and the handoffs are not transparent to the team."
    235 = "The team has access to cross-functional skills and cross-component skills even though these skills are not embedded in the teams. These skills are readily available to assist the team in delivering the functionality. For example the team can access quality assurance engineers as needed without them being embedded in the team."
    236 = "All the handoffs are identified and mapped on the team deployment/release process document."
    237 = "The concept of triads is established in the team and is respected for features to be considered as done. Business, Technical and Quality are committed to execute the strategy together."
    238 = "Team member allocations are 100% to the team."
    239 = "Cross-functional and cross-component skills are embedded in the teams and all the requisite skills required to deliver end-to-end customer value are embedded in the teams. There is no external team dependency to deliver the goals and vision of the product."
    240 = "An onboarding plan for new members exists and this is maintained and reviewed often. The onboarding plan is documented and accessible on the document management application."
    241 = "The triad (business, technology and quality) is fully functional in the teams and the product owners within the teams are empowered to make the decisions pertaining to the product roadmap with need of external consultation or validation."
    242 = "The product goals are clearly prioritized and recorded on the team document management system and the top priority goals are visual on the team walls."
    243 = "Team members can perform more than one function. For example, the developers can serve as quality engineers."
    244 = "Team has autonomy and is self-organized. The team composition is well understood by all team members and stakeholders and the roles are responsibilities are clearly defined and understood."
    245 = "Kaizen (constant improvement) is performed by the team every iteration. The team collects metrics, which are radiated on the wall and regularly holds retrospectives with action items visible on the team wall towards continuous improvement."
    246 = "The team is completely autonomous and drives all decision making. Consultations are used to get information that can feed into their decision making but the decision making is entirely up to the teams."
    247 = "A majority of the team is comprised of individuals who can play dual roles and there is collective shared responsibility by the team members. Failure or success is team failure or success."
    248 = "The team is in charge of the end-to-end delivery of the product from infrastructure, to configuration, to application deployment to instrumentation and alerting."
    249 = "Team performs restrospectives every iteration to drive continuous improvements."
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 5).Value = $values[$row]
}

$heights = @{
    231 = 43.2
    232 = 28.8
    233 = 28.8
    234 = 57.6
    235 = 72
    236 = 28.8
    239 = 57.6
    240 = 43.2
    241 = 57.6
    242 = 43.2
    245 = 57.6
    246 = 43.2
    247 = 43.2
    248 = 43.2
}

foreach ($row in $heights.Keys) {
    $ws.Rows.Item($row).RowHeight = $heights[$row]
}

$ws.Range("E250").Select()
$ws.Application.ActiveWindow.ScrollRow = 242
